$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add Panel Accessories ("MX-DPBX" / "MX-BBX") rows to the Hungary,
# Spain and Turkey market sheets, inserted just above the trailing
# "Wg" / "Accessories" marker rows.
# ---------------------------------------------------------------------

# --- Hungary --------------------------------------------------------
$ws = $wb.Worksheets.Item("Hungary")
$ws.Rows("11:12").Insert()
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)
$ws.Range("A11").Value = "MX-DPBX"
$ws.Range("A12").Value = "MX-BBX"

# --- Spain ------------------------------------------------------------
$ws = $wb.Worksheets.Item("Spain")
$ws.Rows("11:12").Insert()
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)
$ws.Range("A11").Value = "MX-DPBX"
$ws.Range("A12").Value = "MX-BBX"

# --- Turkey -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Turkey")
$ws.Rows("10:11").Insert()
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$ws.Range("A10").Value = "MX-DPBX"
$ws.Range("A11").Value = "MX-BBX"

# ---------------------------------------------------------------------
# Update sheet selections to reflect where the user left the cursor
# after the edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Hungary").Range("A11:A12").Select()
$wb.Worksheets.Item("Poland").Range("A9:A10").Select()
$wb.Worksheets.Item("Turkey").Range("A10").Select()

# Spain becomes the active sheet/tab with the cursor on E4.
$ws = $wb.Worksheets.Item("Spain")
$ws.Activate()
$ws.Range("E4").Select()
